$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 21.857143
$ws.Range("I11").Value = 21.857143
$ws.Range("K11").Value = 21.857143
$ws.Range("M11").Value = 118.142857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 463.6875
$ws.Range("I33").Value = 454.6
$ws.Range("K33").Value = 454.6
$ws.Range("M33").Value = -225.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3118.1
$ws.Range("I38").Value = 123.818184
$ws.Range("J38").Value = 6777.778
$ws.Range("K38").Value = 371.454552
$ws.Range("L38").Value = 20333.334
$ws.Range("M38").Value = 0.5454479999999648
$ws.Range("N38").Value = -21077.334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1997.7894
$ws.Range("I137").Value = 1553.2759
$ws.Range("K137").Value = 4659.8277
$ws.Range("M137").Value = -2109.8277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2817.1667
$ws.Range("I74").Value = 2184.5715
$ws.Range("J74").Value = 3702.8
$ws.Range("K74").Value = 2184.5715
$ws.Range("L74").Value = 3702.8
$ws.Range("M74").Value = -1310.5715
$ws.Range("N74").Value = -5450.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2817.1667
$ws.Range("I77").Value = 2184.5715
$ws.Range("J77").Value = 3702.8
$ws.Range("K77").Value = 10922.8575
$ws.Range("L77").Value = 18514
$ws.Range("M77").Value = -6554.8575
$ws.Range("N77").Value = -27250

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 11906311
$ws.Range("I102").Value = 15152577
$ws.Range("K102").Value = 15152577
$ws.Range("M102").Value = -15150955

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 32757.8
$ws.Range("J108").Value = 32757.8
$ws.Range("L108").Value = 32757.8
$ws.Range("N108").Value = -40437.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3323.75
$ws.Range("I122").Value = 3145.524
$ws.Range("J122").Value = 4571.3335
$ws.Range("K122").Value = 9436.572
$ws.Range("L122").Value = 13714.0005
$ws.Range("M122").Value = -6986.572
$ws.Range("N122").Value = -18614.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2776.318
$ws.Range("I132").Value = 2109.9375
$ws.Range("K132").Value = 6329.8125
$ws.Range("M132").Value = -3799.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 10738
$ws.Range("I29").Value = 1458
$ws.Range("K29").Value = 1458
$ws.Range("M29").Value = -1169

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 30166
$ws.Range("I58").Value = 10000
$ws.Range("J58").Value = 40249
$ws.Range("K58").Value = 10000
$ws.Range("L58").Value = 40249
$ws.Range("N58").Value = -40837
$ws.Range("M58").Value = -9706

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8621339
$ws.Range("I94").Value = 13889414
$ws.Range("J94").Value = 851.8182
$ws.Range("K94").Value = 13889414
$ws.Range("L94").Value = 851.8182
$ws.Range("M94").Value = -13888963
$ws.Range("N94").Value = -1753.8182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 37038000
$ws.Range("I99").Value = 52632388
$ws.Range("K99").Value = 52632388
$ws.Range("M99").Value = -52630890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 200021150
$ws.Range("I105").Value = 200021150
$ws.Range("K105").Value = 200021150
$ws.Range("M105").Value = -200019403

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H116").Value = 70000
$ws.Range("J116").Value = 70000
$ws.Range("L116").Value = 70000
$ws.Range("N116").Value = -79178

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5529.478
$ws.Range("I134").Value = 1196.9445
$ws.Range("K134").Value = 3590.8335
$ws.Range("M134").Value = -1055.8335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1505.8846
$ws.Range("I31").Value = 1457.2084
$ws.Range("J31").Value = 2090
$ws.Range("K31").Value = 1457.2084
$ws.Range("L31").Value = 2090
$ws.Range("M31").Value = -1162.2084
$ws.Range("N31").Value = -2680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1505.8846
$ws.Range("I34").Value = 1457.2084
$ws.Range("J34").Value = 2090
$ws.Range("K34").Value = 1457.2084
$ws.Range("L34").Value = 2090
$ws.Range("M34").Value = -1255.2084
$ws.Range("N34").Value = -2494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2436.087
$ws.Range("I132").Value = 2027
$ws.Range("J132").Value = 3371.1428
$ws.Range("K132").Value = 6081
$ws.Range("L132").Value = 10113.4284
$ws.Range("M132").Value = -3551
$ws.Range("N132").Value = -15173.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1500.5385
$ws.Range("J34").Value = 2468.5715
$ws.Range("L34").Value = 7405.7145
$ws.Range("N34").Value = -7573.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 780.5263
$ws.Range("I41").Value = 480
$ws.Range("J41").Value = 797.2222
$ws.Range("K41").Value = 1440
$ws.Range("L41").Value = 2391.6666
$ws.Range("M41").Value = -1102
$ws.Range("N41").Value = -3067.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2362.0625
$ws.Range("J55").Value = 2806.8462
$ws.Range("L55").Value = 8420.5386
$ws.Range("N55").Value = -8774.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 406.55554
$ws.Range("I114").Value = 248.8
$ws.Range("J114").Value = 603.75
$ws.Range("K114").Value = 746.4000000000001
$ws.Range("L114").Value = 1811.25
$ws.Range("M114").Value = 2507.6
$ws.Range("N114").Value = -8319.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1227.4615
$ws.Range("I122").Value = 596.4286
$ws.Range("J122").Value = 1459.9474
$ws.Range("K122").Value = 5367.8574
$ws.Range("L122").Value = 13139.5266
$ws.Range("M122").Value = -2917.8574
$ws.Range("N122").Value = -18039.5266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 55950
$ws.Range("J124").Value = 55950
$ws.Range("L124").Value = 55950
$ws.Range("N124").Value = -65770

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3424.7
$ws.Range("I132").Value = 3214.318
$ws.Range("J132").Value = 4003.25
$ws.Range("K132").Value = 9642.954000000002
$ws.Range("L132").Value = 12009.75
$ws.Range("M132").Value = -7112.954000000002
$ws.Range("N132").Value = -17069.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2757
$ws.Range("I7").Value = 2293.3333
$ws.Range("J7").Value = 3452.5
$ws.Range("K7").Value = 2293.3333
$ws.Range("L7").Value = 3452.5
$ws.Range("M7").Value = -2181.3333
$ws.Range("N7").Value = -3676.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1144.4286
$ws.Range("I22").Value = 1002.75
$ws.Range("J22").Value = 1333.3334
$ws.Range("K22").Value = 1002.75
$ws.Range("L22").Value = 1333.3334
$ws.Range("M22").Value = -707.75
$ws.Range("N22").Value = -1923.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1144.4286
$ws.Range("I27").Value = 1002.75
$ws.Range("J27").Value = 1333.3334
$ws.Range("K27").Value = 1002.75
$ws.Range("L27").Value = 1333.3334
$ws.Range("M27").Value = -895.75
$ws.Range("N27").Value = -1547.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1258.75
$ws.Range("I61").Value = 1150.7
$ws.Range("J61").Value = 1438.8334
$ws.Range("K61").Value = 1150.7
$ws.Range("L61").Value = 1438.8334
$ws.Range("M61").Value = -948.7
$ws.Range("N61").Value = -1842.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 930.53845
$ws.Range("I93").Value = 966.8889
$ws.Range("J93").Value = 848.75
$ws.Range("K93").Value = 966.8889
$ws.Range("L93").Value = 848.75
$ws.Range("M93").Value = 281.1111
$ws.Range("N93").Value = -3344.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1258.75
$ws.Range("I113").Value = 1150.7
$ws.Range("J113").Value = 1438.8334
$ws.Range("K113").Value = 1150.7
$ws.Range("L113").Value = 1438.8334
$ws.Range("M113").Value = 1019.3
$ws.Range("N113").Value = -5778.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 27779506
$ws.Range("I122").Value = 41668292
$ws.Range("J122").Value = 1935
$ws.Range("K122").Value = 125004876
$ws.Range("L122").Value = 5805
$ws.Range("M122").Value = -125002426
$ws.Range("N122").Value = -10705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2757
$ws.Range("I126").Value = 2293.3333
$ws.Range("J126").Value = 3452.5
$ws.Range("K126").Value = 6879.999899999999
$ws.Range("L126").Value = 10357.5
$ws.Range("M126").Value = -4409.999899999999
$ws.Range("N126").Value = -15297.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 44200
$ws.Range("J130").Value = 44200
$ws.Range("L130").Value = 44200
$ws.Range("N130").Value = -54240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 6886.6665
$ws.Range("J15").Value = 6886.6665
$ws.Range("L15").Value = 6886.6665
$ws.Range("N15").Value = -7462.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 39883.1
$ws.Range("J109").Value = 38721
$ws.Range("L109").Value = 38721
$ws.Range("N109").Value = -41495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9261490
$ws.Range("I122").Value = 11907006
$ws.Range("J122").Value = 2185
$ws.Range("K122").Value = 35721018
$ws.Range("L122").Value = 6555
$ws.Range("M122").Value = -35718568
$ws.Range("N122").Value = -11455

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 40000708
$ws.Range("I126").Value = 47619730
$ws.Range("J126").Value = 852.5
$ws.Range("K126").Value = 142859190
$ws.Range("L126").Value = 2557.5
$ws.Range("M126").Value = -142856720
$ws.Range("N126").Value = -7497.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2089.1052
$ws.Range("I132").Value = 1657.8387
$ws.Range("K132").Value = 4973.5161
$ws.Range("M132").Value = -2443.5161

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 72857.5
$ws.Range("J135").Value = 72857.5
$ws.Range("L135").Value = 72857.5
$ws.Range("N135").Value = -82997.5
